$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 969.06665
$ws.Range("I40").Value = 899.4
$ws.Range("J40").Value = 1003.9
$ws.Range("K40").Value = 899.4
$ws.Range("L40").Value = 1003.9
$ws.Range("M40").Value = -724.4
$ws.Range("N40").Value = -1353.9

# Row 138
$ws.Range("H138").Value = 1775.99
$ws.Range("J138").Value = 2190.5652
$ws.Range("L138").Value = 6571.6956
$ws.Range("N138").Value = -16851.6956


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 14986
$ws.Range("I2").Value = 1176.6
$ws.Range("J2").Value = 49509.5
$ws.Range("K2").Value = 1176.6
$ws.Range("L2").Value = 49509.5
$ws.Range("M2").Value = -1063.6
$ws.Range("N2").Value = -49735.5

# Row 45
$ws.Range("H45").Value = 1043.7858
$ws.Range("I45").Value = 905
$ws.Range("J45").Value = 1182.5714
$ws.Range("K45").Value = 905
$ws.Range("L45").Value = 1182.5714
$ws.Range("M45").Value = -528
$ws.Range("N45").Value = -1936.5714

# Row 116
$ws.Range("H116").Value = 14986
$ws.Range("I116").Value = 1176.6
$ws.Range("J116").Value = 49509.5
$ws.Range("K116").Value = 1176.6
$ws.Range("L116").Value = 49509.5
$ws.Range("M116").Value = 1117.4
$ws.Range("N116").Value = -54097.5


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 14986
$ws.Range("I3").Value = 1176.6
$ws.Range("J3").Value = 49509.5
$ws.Range("K3").Value = 1176.6
$ws.Range("L3").Value = 49509.5
$ws.Range("M3").Value = -1062.6
$ws.Range("N3").Value = -49737.5

# Row 86
$ws.Range("H86").Value = 1001715.3
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2334669
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2334669
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -2336915

# Row 89
$ws.Range("H89").Value = 1001715.3
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2334669
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 11673345
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -11684577

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1007.3333
$ws.Range("I16").Value = 1252.5
$ws.Range("J16").Value = 811.2
$ws.Range("K16").Value = 1252.5
$ws.Range("L16").Value = 811.2
$ws.Range("M16").Value = -965.5
$ws.Range("N16").Value = -1385.2

# Row 86
$ws.Range("H86").Value = 38464290
$ws.Range("I86").Value = 52634430
$ws.Range("J86").Value = 2457.1428
$ws.Range("K86").Value = 52634430
$ws.Range("L86").Value = 2457.1428
$ws.Range("M86").Value = -52633307
$ws.Range("N86").Value = -4703.1428

# Row 89
$ws.Range("H89").Value = 38464290
$ws.Range("I89").Value = 52634430
$ws.Range("J89").Value = 2457.1428
$ws.Range("K89").Value = 263172150
$ws.Range("L89").Value = 12285.714
$ws.Range("M89").Value = -263166534
$ws.Range("N89").Value = -23517.714

# Row 113
$ws.Range("H113").Value = 1007.3333
$ws.Range("I113").Value = 1252.5
$ws.Range("J113").Value = 811.2
$ws.Range("K113").Value = 1252.5
$ws.Range("L113").Value = 811.2
$ws.Range("M113").Value = 917.5
$ws.Range("N113").Value = -5151.2


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 29419788
$ws.Range("I68").Value = 17067.166
$ws.Range("J68").Value = 45457636
$ws.Range("K68").Value = 51201.49800000001
$ws.Range("L68").Value = 136372908
$ws.Range("M68").Value = -50390.49800000001
$ws.Range("N68").Value = -136374530

# Row 71
$ws.Range("H71").Value = 29419788
$ws.Range("I71").Value = 17067.166
$ws.Range("J71").Value = 45457636
$ws.Range("K71").Value = 153604.494
$ws.Range("L71").Value = 409118724
$ws.Range("M71").Value = -149548.494
$ws.Range("N71").Value = -409126836

# Row 113
$ws.Range("H113").Value = 691.375
$ws.Range("I113").Value = 458.33334
$ws.Range("K113").Value = 1375.00002
$ws.Range("M113").Value = 794.9999800000001

# Row 137
$ws.Range("H137").Value = 7947188
$ws.Range("I137").Value = 16667881
$ws.Range("J137").Value = 4458910.5
$ws.Range("K137").Value = 50003643
$ws.Range("L137").Value = 13376731.5
$ws.Range("M137").Value = -49998543
$ws.Range("N137").Value = -13386931.5

# Row 138
$ws.Range("H138").Value = 8774650
$ws.Range("I138").Value = 1238.2354
$ws.Range("K138").Value = 3714.7062
$ws.Range("M138").Value = 1425.2938


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15970

# Row 52
$ws.Range("H52").Value = 33874.5
$ws.Range("I52").Value = 18488.334
$ws.Range("K52").Value = 18488.334
$ws.Range("M52").Value = -18229.334

# Row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 29866.666
$ws.Range("J42").Value = 29866.666
$ws.Range("L42").Value = 29866.666
$ws.Range("N42").Value = -30992.666

# Row 46
$ws.Range("H46").Value = 786
$ws.Range("I46").Value = 650
$ws.Range("K46").Value = 650
$ws.Range("M46").Value = -462

# Row 49
$ws.Range("H49").Value = 29866.666
$ws.Range("J49").Value = 29866.666
$ws.Range("L49").Value = 29866.666
$ws.Range("N49").Value = -30160.666

# Row 50
$ws.Range("H50").Value = 34800
$ws.Range("J50").Value = 34800
$ws.Range("L50").Value = 34800
$ws.Range("M50").Value = -36074

# Row 56
$ws.Range("H56").Value = 36264.25
$ws.Range("I56").Value = 5000
$ws.Range("J56").Value = 46685.668
$ws.Range("K56").Value = 5000
$ws.Range("L56").Value = 46685.668
$ws.Range("M56").Value = -4309
$ws.Range("N56").Value = -48067.668


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 38
$ws.Range("H38").Value = 39265.5
$ws.Range("I38").Value = 3000
$ws.Range("J38").Value = 51354
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 51354
$ws.Range("M38").Value = -2527
$ws.Range("N38").Value = -52300

# Row 54
$ws.Range("H54").Value = 29994.25
$ws.Range("J54").Value = 29994.25
$ws.Range("L54").Value = 29994.25
$ws.Range("N54").Value = -31034.25

# Row 55
$ws.Range("H55").Value = 3725
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3725
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 3725
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -4279

# Row 61
$ws.Range("H61").Value = 14150
$ws.Range("J61").Value = 14150
$ws.Range("L61").Value = 14150
$ws.Range("N61").Value = -14734

